$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write header values in the same order the original author must have used,
# so the shared-strings table comes out in the same index order:
# 0 = "Total de Ventas", 1 = "Monto de Facturación", 2 = "Tipo de Cliente"
$ws.Range("B1").Value = "Total de Ventas"
$ws.Range("C1").Value = "Monto de Facturación"
$ws.Range("A1").Value = "Tipo de Cliente"

# Build the header style once on A1 (bold, 12pt, accent fill, vertically
# centered), then replicate it onto B1:C1 via a format-only paste so every
# header cell resolves to the same single cellXf.
$hdr = $ws.Range("A1")
$hdr.Font.Bold = $true
$hdr.Font.Size = 12
$hdr.Interior.Pattern = 1
$hdr.Interior.ThemeColor = 5
$hdr.Interior.TintAndShade = 0.59999389629810485
$hdr.VerticalAlignment = -4108

$hdr.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths / row height for the header row.
$ws.Columns.Item(1).ColumnWidth = 18.276041666666668
$ws.Columns.Item(2).ColumnWidth = 18.944010416666668
$ws.Columns.Item(3).ColumnWidth = 23.944010416666668
$ws.Rows.Item(1).RowHeight = 22.8

# Match the saved selection state.
$ws.Range("B8").Select()
